$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 163, shifting existing rows 163:178 down to 164:179
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new data record
$ws.Cells.Item(163, 1).Value = 4
$ws.Cells.Item(163, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(163, 3).Value = "Los Lagos"
$ws.Cells.Item(163, 4).Value = 44578
$ws.Cells.Item(163, 5).Value = 10
$ws.Cells.Item(163, 6).Value = 100112032
$ws.Cells.Item(163, 7).Value = "Zapallo italiano"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 150
$ws.Cells.Item(163, 11).Value = 12000
$ws.Cells.Item(163, 12).Value = 12000
$ws.Cells.Item(163, 13).Value = 12000
$ws.Cells.Item(163, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 200
$ws.Cells.Item(163, 17).Value = 60
$ws.Cells.Item(163, 18).Value = "Hortaliza"

# Apply the date number format to column D of the new row, matching the rest of the column
$ws.Cells.Item(163, 4).NumberFormat = $ws.Cells.Item(164, 4).NumberFormat
